$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.554.47'
$ws.Range('E2').Value = '  +1.60%  '
$ws.Range('D3').Value = '3.024.66'
$ws.Range('E3').Value = '  +3.59%  '
$ws.Range('E4').Value = '  +0.07%  '
$__style = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '380.06'
$ws.Range('D5').Style = $__style
$ws.Range('E5').Value = '  +1.36%  '
$__style = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '103.29'
$ws.Range('D6').Style = $__style
$ws.Range('E6').Value = '  +3.78%  '
$__style = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.545'
$ws.Range('D7').Style = $__style
$ws.Range('E7').Value = '  +1.99%  '
$ws.Range('E8').Value = '  +0.00%  '
$__style = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.595'
$ws.Range('D9').Style = $__style
$ws.Range('E9').Value = '  +3.78%  '
$__style = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.84'
$ws.Range('D10').Style = $__style
$ws.Range('E10').Value = '  +3.53%  '
$ws.Range('E11').Value = '  -0.03%  '
$__style = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0858'
$ws.Range('D12').Style = $__style
$ws.Range('E12').Value = '  +1.53%  '
$ws.Range('D13').Value = '3.503.79'
$ws.Range('E13').Value = '  +3.67%  '
$__style = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '18.56'
$ws.Range('D14').Style = $__style
$ws.Range('E14').Value = '  +3.37%  '
$__style = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.77'
$ws.Range('D15').Style = $__style
$ws.Range('E15').Value = '  +2.27%  '
$ws.Range('D16').Value = '3.037.29'
$ws.Range('E16').Value = '  +3.81%  '
$__style = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.993'
$ws.Range('D17').Style = $__style
$ws.Range('E17').Value = '  +0.37%  '
$ws.Range('E18').Value = '  -10.45%  '
$ws.Range('D19').Value = '51.646.90'
$ws.Range('E19').Value = '  +1.88%  '
$__style = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.06'
$ws.Range('D20').Style = $__style
$ws.Range('E20').Value = '  +2.21%  '
$ws.Range('E21').Value = '  +2.14%  '
$ws.Range('E22').Value = '  +2.42%  '
$__style = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.25'
$ws.Range('D23').Style = $__style
$ws.Range('E23').Value = '  +1.52%  '
$__style = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '268.51'
$ws.Range('D24').Style = $__style
$ws.Range('E24').Value = '  +1.12%  '
$ws.Range('E25').Value = '  +0.81%  '
$__style = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.30'
$ws.Range('D26').Style = $__style
$ws.Range('E26').Value = '  +5.96%  '
$__style = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.47'
$ws.Range('D27').Style = $__style
$ws.Range('E27').Value = '  +5.79%  '
$ws.Range('E28').Value = '  +6.23%  '
$ws.Range('E29').Value = '  -0.13%  '
$__style = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '26.21'
$ws.Range('D30').Style = $__style
$ws.Range('E30').Value = '  +3.55%  '
$ws.Range('E31').Value = '  +2.14%  '
$__style = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '10.33'
$ws.Range('D32').Style = $__style
$ws.Range('E32').Value = '  +4.18%  '
$ws.Range('B33').Value = 'InjectiveProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$__style = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '34.32'
$ws.Range('D33').Style = $__style
$ws.Range('E33').Value = '  +3.39%  '
$ws.Range('B34').Value = 'OKB'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$__style = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '51.08'
$ws.Range('D34').Style = $__style
$ws.Range('E34').Value = '  +1.56%  '
$__style = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.05'
$ws.Range('D35').Style = $__style
$ws.Range('E35').Value = '  +0.60%  '
$__style = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0449'
$ws.Range('D36').Style = $__style
$ws.Range('E36').Value = '  +5.13%  '
$ws.Range('E37').Value = '  +0.08%  '
$__style = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.29'
$ws.Range('D38').Style = $__style
$ws.Range('E38').Value = '  +7.85%  '
$__style = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '17.31'
$ws.Range('D39').Style = $__style
$ws.Range('E39').Value = '  +6.42%  '
$__style = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.61'
$ws.Range('D40').Style = $__style
$ws.Range('E40').Value = '  +8.57%  '
$ws.Range('B41').Value = 'TheGraph'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$__style = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.283'
$ws.Range('D41').Style = $__style
$ws.Range('E41').Value = '  +10.22%  '
$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$__style = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.87'
$ws.Range('D42').Style = $__style
$ws.Range('E42').Value = '  +4.79%  '
$ws.Range('E43').Value = '  +1.34%  '
$__style = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '126.80'
$ws.Range('D44').Style = $__style
$ws.Range('E44').Value = '  +3.14%  '
$__style = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.77'
$ws.Range('D45').Style = $__style
$ws.Range('E45').Value = '  +13.39%  '
$__style = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '22.08'
$ws.Range('D46').Style = $__style
$ws.Range('E46').Value = '  +6.28%  '
$ws.Range('E47').Value = '  +0.80%  '
$ws.Range('E48').Value = '  +2.45%  '
$ws.Range('D49').Value = '2.038.96'
$ws.Range('E49').Value = '  +2.36%  '
$ws.Range('D50').Value = '3.328.48'
$ws.Range('E50').Value = '  +3.70%  '
$ws.Range('E51').Value = '  +3.13%  '
